# Apply cell updates from the Jan 3, 2023 GitHub Actions symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range('D2').Value = '''245.96'
$ws.Range('E2').Value = '''-0.30%'

# Row 3
$ws.Range('D3').Value = '''29.76'
$ws.Range('E3').Value = '''-1.26%'

# Row 4
$ws.Range('D4').Value = '''5.152'
$ws.Range('E4').Value = '''-0.47%'

# Row 5
$ws.Range('D5').Value = '''0.05772'
$ws.Range('E5').Value = '''0.28%'

# Row 6
$ws.Range('E6').Value = '''1.00%'

# Row 7
$ws.Range('E7').Value = '''4.63%'

# Row 8
$ws.Range('D8').Value = '''0.8536'
$ws.Range('E8').Value = '''-0.33%'

# Row 9
$ws.Range('D9').Value = '''0.8542'
$ws.Range('E9').Value = '''-3.06%'

# Row 10
$ws.Range('D10').Value = '''0.1377'
$ws.Range('E10').Value = '''0.72%'

# Row 11
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '''0.03375'
$ws.Range('E11').Value = '''2.24%'

# Row 12
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.07082'
$ws.Range('E12').Value = '''1.42%'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03262'
$ws.Range('E13').Value = '''11.40%'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09365'
$ws.Range('E14').Value = '''-0.28%'

# Row 15
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001530'
$ws.Range('E15').Value = '''1.31%'

# Row 16
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D16').Value = '''0.0005963'
$ws.Range('E16').Value = '''-0.31%'

# Row 17
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = '''0.006036'
$ws.Range('E17').Value = '''-0.47%'

# Row 18
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = '''3.512'
$ws.Range('E18').Value = '''0.09%'

# Row 19
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').Value = '''2.222'
$ws.Range('E19').Value = '''1.96%'

# Row 20
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D20').Value = '''0.3158'
$ws.Range('E20').Value = '''0.35%'

# Row 22
$ws.Range('D22').Value = '''3.491'
$ws.Range('E22').Value = '''-3.66%'

# Row 23
$ws.Range('B23').Value = 'ZBToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D23').Value = '''0.1411'
$ws.Range('E23').Value = '''2.45%'

# Row 24
$ws.Range('B24').Value = 'CoinExToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D24').Value = '''0.04123'
$ws.Range('E24').Value = '''-1.07%'

# Row 25
$ws.Range('D25').Value = '''0.001229'
$ws.Range('E25').Value = '''1.26%'

# Row 26
$ws.Range('D26').Value = '''0.004139'
$ws.Range('E26').Value = '''-8.07%'

# Row 27
$ws.Range('D27').Value = '''0.0001200'
$ws.Range('E27').Value = '''1.86%'

# Row 28
$ws.Range('D28').Value = '''0.0001449'
$ws.Range('E28').Value = '''5.19%'

# Row 40
$ws.Range('D40').Value = '''0.03750'
$ws.Range('E40').Value = '''-0.87%'

# Row 41
$ws.Range('D41').Value = '''0.005662'
$ws.Range('E41').Value = '''61.29%'

# Row 42
$ws.Range('D42').Value = '''0.1069'
$ws.Range('E42').Value = '''-0.01%'

# Row 43
$ws.Range('D43').Value = '''0.002301'
$ws.Range('E43').Value = '''-11.05%'

# Row 44
$ws.Range('D44').Value = '''0.008494'
$ws.Range('E44').Value = '''-15.19%'

# Row 45
$ws.Range('D45').Value = '''0.00005430'
$ws.Range('E45').Value = '''6.43%'

# Row 46
$ws.Range('E46').Value = '''0.19%'

# Row 47
$ws.Range('E47').Value = '''-20.07%'

# Row 48
$ws.Range('D48').Value = '''0.002214'
$ws.Range('E48').Value = '''-18.56%'

# Row 49
$ws.Range('E49').Value = '''0.19%'

# Row 50
$ws.Range('E50').Value = '''0.19%'
